# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet right before the "总计" (total) sheet,
#    populated with the Q1-2022 fund-holding detail rows.
# 2. Prepend a "2022-Q1" summary row to the "总计" sheet (pushing the
#    existing quarters down and renumbering the index column).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, inserted immediately before "总计"
# ---------------------------------------------------------------------
$totalBeforeInsert = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalBeforeInsert)
$q1.Name = "2022-Q1"

# Sheet references returned by Worksheets.Item() track the tab POSITION,
# not the sheet identity, so after Add() shifts "总计" one slot to the
# right we must re-fetch it by name to keep pointing at the right sheet.
$total = $wb.Worksheets.Item("总计")

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $q1.Range($cols[$i] + "1").Value = $headers[$i]
}

# header row look: bold, centered, top-aligned, thin box border (matches
# the style already used for header / index cells on the other sheets)
$headerRange = $q1.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1

# Fund holding rows: code, name, scale, stock position, position ratio,
# held market value (billion), position rank
$rows = @(
    @("164705", "汇添富恒生指数（QDII-LOF）A", "2.96", "92.23", "7.67", "0.2270", 1),
    @("160924", "大成恒生指数（QDII-LOF）", "0.89", "93.20", "7.89", "0.0702", 1),
    @("008253", "华宝致远混合（QDII）A", "0.70", "85.00", "4.67", "0.0327", 5),
    @("010789", "汇添富恒生指数（QDII-LOF）C", "0.37", "92.23", "7.67", "0.0284", 1),
    @("008254", "华宝致远混合（QDII）C", "0.12", "85.00", "4.67", "0.0056", 5)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $data = $rows[$r]

    $q1.Range("A" + $rowNum).Value = $r
    $q1.Range("A" + $rowNum).Font.Bold = $true
    $q1.Range("A" + $rowNum).HorizontalAlignment = -4108
    $q1.Range("A" + $rowNum).VerticalAlignment = -4160
    $q1.Range("A" + $rowNum).Borders.LineStyle = 1

    # text-valued columns: code / name / scale / stock position / ratio / value
    $q1.Range("B" + $rowNum).NumberFormat = "@"
    $q1.Range("B" + $rowNum).Value = $data[0]

    $q1.Range("C" + $rowNum).Value = $data[1]

    $q1.Range("D" + $rowNum).NumberFormat = "@"
    $q1.Range("D" + $rowNum).Value = $data[2]

    $q1.Range("E" + $rowNum).NumberFormat = "@"
    $q1.Range("E" + $rowNum).Value = $data[3]

    $q1.Range("F" + $rowNum).NumberFormat = "@"
    $q1.Range("F" + $rowNum).Value = $data[4]

    $q1.Range("G" + $rowNum).NumberFormat = "@"
    $q1.Range("G" + $rowNum).Value = $data[5]

    # position rank is a real number
    $q1.Range("H" + $rowNum).Value = $data[6]
}

# ---------------------------------------------------------------------
# 2) "总计" sheet: insert a new row for 2022-Q1 ahead of the existing data
# ---------------------------------------------------------------------
$total.Rows(2).Insert()
$total.Range("A2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.36

# renumber the index column for the quarters that got pushed down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# restore the index-column style (bold/border/center-top) on the new A2,
# copying it from a cell that already carries it
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
